# "ran model for jan 7"
# Fill in "Beat Vegas?" predictions for the already-listed Jan 6 games
# (rows 7-17), correct the Jan 6 Warriors @ Clippers prediction (row 17),
# and append the newly-modeled Jan 7 games (rows 18-22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in "Beat Vegas?" (column G) for rows that already had a
#     Spread/Predicted Spread but were missing the outcome column ---
$beatVegas = @{
    7  = "Yes"
    8  = "Yes"
    9  = "Yes"
    10 = "Yes"
    11 = "No"
    12 = "Yes"
    13 = "Yes"
    14 = "Yes"
    15 = "Yes"
    16 = "Yes"
}
foreach ($row in $beatVegas.Keys) {
    $ws.Cells.Item($row, 7).Value = $beatVegas[$row]
}

# --- Row 17 (GSW vs LAC, Jan 6): predicted spread / difference updated ---
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 6).Value = 8.5
$ws.Cells.Item(17, 7).Value = "No"

# --- New Jan 7 (serial 44203) games modeled ---
$newGames = @(
    @(44203, "BRK", "PHI", 1.5, 2.6, -1.1000000000000001),
    @(44203, "MEM", "CLE", -5, -7.9, 2.9),
    @(44203, "POR", "MIN", -10, -4.4000000000000004, -5.6),
    @(44203, "DEN", "DAL", -2.5, -11.4, 8.9),
    @(44203, "LAL", "SAS", -7, -6.1, -0.90000000000000036)
)

$r = 18
foreach ($game in $newGames) {
    $ws.Cells.Item($r, 1).Value = $game[0]
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy\-mm\-dd"
    $ws.Cells.Item($r, 2).Value = $game[1]
    $ws.Cells.Item($r, 3).Value = $game[2]
    $ws.Cells.Item($r, 4).Value = $game[3]
    $ws.Cells.Item($r, 5).Value = $game[4]
    $ws.Cells.Item($r, 6).Value = $game[5]
    $r++
}

# Restore the user's final selection
$ws.Range("K14").Select()
